$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.49 = 30396.71 pesos`n✅ 30396.71 pesos = 7.45 = 979.4 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate cells N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 133.6
$wsTasas.Range("O10").Value = 4061
$wsTasas.Range("N12").Value = 4080
$wsTasas.Range("O12").Value = 131.46
